$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-09-23 Tuesday"; new = "2025-09-24 Wednesday"},
    @{old = "83×31=";             new = "92×54="},
    @{old = "18×72=";             new = "26×54="},
    @{old = "70×96=";             new = "50×66="},
    @{old = "48×68=";             new = "47×12="},
    @{old = "30×62=";             new = "52×60="},
    @{old = "47×65=";             new = "13×74="},
    @{old = "97×48=";             new = "95×45="},
    @{old = "28×55=";             new = "82×48="},
    @{old = "52×73=";             new = "69×61="},
    @{old = "89×24=";             new = "43×19="},
    @{old = "15×52=";             new = "66×49="},
    @{old = "52×13=";             new = "93×59="},
    @{old = "81×31=";             new = "62×26="},
    @{old = "86×60=";             new = "57×25="},
    @{old = "13×19=";             new = "31×85="},
    @{old = "54×57=";             new = "99×35="},
    @{old = "41×37=";             new = "31×85="},
    @{old = "79×88=";             new = "90×69="},
    @{old = "68×57=";             new = "86×55="},
    @{old = "97×34=";             new = "26×75="},
    @{old = "60×66=";             new = "99×65="},
    @{old = "41×71=";             new = "61×66="},
    @{old = "17×57=";             new = "65×37="},
    @{old = "42×87=";             new = "93×12="},
    @{old = "60×44=";             new = "74×36="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
